$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C260) from 45203 to 45205 (2023-10-04 -> 2023-10-06)
$ws.Range("C2:C260").Value = 45205
